$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Numeric-looking D values are forced to text (matching the sheet's existing
# text-based price formatting) via a temporary "@" number format, then the
# format is cleared again so no residual style is left on the cell.

$ws.Range("D2").Value = "27.046.35"
$ws.Range("E2").Value = "  -3.15%  "

$ws.Range("D3").Value = "1.742.19"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9975"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.60"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9969"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4947"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3521"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.78"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07271"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.061"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9967"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.08"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.906"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").Value = "1.736.00"
$ws.Range("E15").Value = "  -1.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.843"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.29"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.44%  "

$ws.Range("E18").Value = "  -1.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06385"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9971"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("E21").Value = "  -1.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.742"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("D23").Value = "27.096.10"
$ws.Range("E23").Value = "  -3.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.98"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.057"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.71"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.07"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").Value = "1.931.20"
$ws.Range("E28").Value = "  -2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.092"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.47"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.055"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09373"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.578"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.411"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05941"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02194"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.97%  "

$ws.Range("E37").Value = "  -5.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.433"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1995"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.766"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6035"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9965"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.115"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.451"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.89"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.575"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5646"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.849"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06678"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.91%  "

$ws.Range("E51").Value = "  -2.62%  "
